# Resolve missing idAttribute errors
#
# The ISA-TAB entity names ("Investigation", "Study_Section", "Study_Node",
# "Assay_Node", "Assay_Microarray", "Assay_Electrophoresis", "Assay_MS",
# "Assay_NMR") are renamed with an "ISATAB_" prefix so that the attribute
# sheet's "entity" / "refEntity" columns (and the entities sheet's "extends"
# column) line up with the renamed entities and the idAttribute lookup no
# longer fails.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("entities")
$ws2 = $wb.Worksheets.Item("attributes")

# old name -> new (ISATAB_-prefixed) name, applied in both sheets
$renames = [ordered]@{
    "Investigation"          = "ISATAB_Investigation"
    "Study_Section"          = "ISATAB_Study_Section"
    "Study_Node"             = "ISATAB_Study_Node"
    "Assay_Node"             = "ISATAB_Assay_Node"
    "Assay_Microarray"       = "ISATAB_Assay_Microarray"
    "Assay_Electrophoresis"  = "ISATAB_Assay_Electrophoresis"
    "Assay_MS"               = "ISATAB_Assay_MS"
    "Assay_NMR"              = "ISATAB_Assay_NMR"
}

$xlWhole = 1

# Apply the renames entity-by-entity (Investigation, Study_Section, ...) so
# new shared-string entries are minted in that same order, matching how the
# workbook's string table grows when an author edits row-by-row.
#   entities!D2:D9   - the "extends" column
#   attributes!C2:C131 - the "entity" column
#   attributes!E2:E131 - the "refEntity" column
foreach ($old in $renames.Keys) {
    $new = $renames[$old]
    $ws1.Range("D2:D9").Replace($old, $new, $xlWhole) | Out-Null
    $ws2.Range("C2:C131").Replace($old, $new, $xlWhole) | Out-Null
    $ws2.Range("E2:E131").Replace($old, $new, $xlWhole) | Out-Null
}

# Column C/E now hold longer strings - widen them to fit, mirroring Excel's
# own "best fit" column behaviour after the edit.
$ws2.Columns.Item(3).ColumnWidth = 25.5
$ws2.Columns.Item(5).ColumnWidth = 18.3

# Reflect the author's final selection / active sheet.
$ws2.Range("A6").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("D12").Select() | Out-Null
